# "chỉnh phần báo cáo" - update the Transaction report table:
#  - remove the "transt0" transaction row
#  - remove the "transt2" transaction row
#  - append a new "transt6" transaction row at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds "transt0" (Rut tien / 11/10/2024 / 42134124 / sfafasdfa) -> delete it.
$ws.Rows(2).Delete()

# After the delete above, the former "transt2" row (Nap tien / 11/10/2024 /
# 241512413 / fsafsasd) has shifted up to row 3 -> delete it as well.
$ws.Rows(3).Delete()

# The remaining data now occupies rows 2-5 (transt1, transt3, transt4, transt5).
# Append the new transaction "transt6" on row 6.
$ws.Range("A6").Value = "transt6"
$ws.Range("B6").Value = "Rút tiền"
$ws.Range("C6").Value = "11/10/2024 12:00:00 AM"

# "4311342" looks like a number, so force it to be stored as text (matching
# how every other amount/id/date column in this sheet is stored) without
# leaving a permanent custom number format behind.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4311342"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "sdfasfdasf"
